$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F29").Value = 32
$ws.Range("G29").Value = 3278.4
$ws.Range("F37").Value = 19
$ws.Range("G37").Value = 875.9
$ws.Range("F40").Value = 29
$ws.Range("G40").Value = 965.7
$ws.Range("B41").Value = 21889.41
$ws.Range("F50").Value = 169
$ws.Range("G50").Value = 32598.41
$ws.Range("F67").Value = 40
$ws.Range("G67").Value = 2230.4
$ws.Range("B77").Value = 129094.84
$ws.Range("F104").Value = 85
$ws.Range("G104").Value = 3383
$ws.Range("B105").Value = 8433.18
$ws.Range("F113").Value = 170
$ws.Range("G113").Value = 7633
$ws.Range("B116").Value = 76967.92999999999
$ws.Range("F119").Value = 7
$ws.Range("G119").Value = 165375.63
$ws.Range("B122").Value = 989002.95
$ws.Range("F126").Value = 24
$ws.Range("G126").Value = 1265.04
$ws.Range("F146").Value = 72
$ws.Range("G146").Value = 1388.16
$ws.Range("B147").Value = 58374.84
$ws.Range("F157").Value = 57
$ws.Range("G157").Value = 2820.36
$ws.Range("F160").Value = 26
$ws.Range("G160").Value = 1100.32
$ws.Range("F165").Value = 30
$ws.Range("G165").Value = 2747.1
$ws.Range("B168").Value = 29252.58
$ws.Range("B179").Value = 57756
$ws.Range("B180").Value = 53925
$ws.Range("F200").Value = 11
$ws.Range("G200").Value = 860.86
$ws.Range("F205").Value = 42
$ws.Range("G205").Value = 3286.92
$ws.Range("F207").Value = 38
$ws.Range("G207").Value = 2973.88
$ws.Range("B210").Value = 30051.89
$ws.Range("F236").Value = 26
$ws.Range("G236").Value = 1089.66
$ws.Range("B237").Value = 14216.31
$ws.Range("F240").Value = 15
$ws.Range("G240").Value = 518.25
$ws.Range("B244").Value = 1898.96
$ws.Range("F255").Value = 72
$ws.Range("G255").Value = 1782.72
$ws.Range("B264").Value = 89844.37
$ws.Range("F272").Value = 7
$ws.Range("G272").Value = 24461.5
$ws.Range("B284").Value = 206372.46
$ws.Range("F287").Value = 46
$ws.Range("G287").Value = 6234.84
$ws.Range("F318").Value = 31
$ws.Range("G318").Value = 7668.78
$ws.Range("B336").Value = 253978.78
$ws.Range("F351").Value = 70
$ws.Range("G351").Value = 1618.4
$ws.Range("F355").Value = 13
$ws.Range("G355").Value = 1649.18
$ws.Range("F362").Value = 43
$ws.Range("G362").Value = 3983.09
$ws.Range("F363").Value = 316
$ws.Range("G363").Value = 12561
$ws.Range("B364").Value = 61985.11
$ws.Range("F368").Value = 3
$ws.Range("G368").Value = 676.71
$ws.Range("B370").Value = 4810.34
$ws.Range("F378").Value = 4
$ws.Range("G378").Value = 85
$ws.Range("F385").Value = 189
$ws.Range("G385").Value = 6822.9
$ws.Range("F396").Value = 11
$ws.Range("G396").Value = 1374.01
$ws.Range("B397").Value = 26577.37
$ws.Range("F401").Value = 66
$ws.Range("G401").Value = 4006.2
$ws.Range("F404").Value = 16
$ws.Range("G404").Value = 2202.72
$ws.Range("F406").Value = 52
$ws.Range("G406").Value = 7814.04
$ws.Range("B407").Value = 16197.4
$ws.Range("F419").Value = 21
$ws.Range("G419").Value = 1402.59
$ws.Range("B421").Value = 31170.85
$ws.Range("F431").Value = 190
$ws.Range("G431").Value = 5650.6
$ws.Range("F435").Value = 169
$ws.Range("G435").Value = 4155.71
$ws.Range("B439").Value = 131915.28
$ws.Range("F461").Value = 454
$ws.Range("G461").Value = 5856.6
$ws.Range("F463").Value = 87
$ws.Range("G463").Value = 4397.85
$ws.Range("F464").Value = 651
$ws.Range("G464").Value = 8397.9
$ws.Range("F465").Value = 170
$ws.Range("G465").Value = 5948.3
$ws.Range("F469").Value = 106
$ws.Range("G469").Value = 2036.26
$ws.Range("B470").Value = 39763.36
$ws.Range("F479").Value = 798
$ws.Range("G479").Value = 10493.7
$ws.Range("F480").Value = 898
$ws.Range("G480").Value = 11503.38
$ws.Range("F481").Value = 572
$ws.Range("G481").Value = 15043.6
$ws.Range("F482").Value = 472
$ws.Range("G482").Value = 9312.559999999999
$ws.Range("F488").Value = 288
$ws.Range("G488").Value = 5604.48
$ws.Range("F489").Value = 1428
$ws.Range("G489").Value = 9396.24
$ws.Range("F491").Value = 1437
$ws.Range("G491").Value = 9326.129999999999
$ws.Range("F492").Value = 575
$ws.Range("G492").Value = 7561.25
$ws.Range("B496").Value = 166223.04
$ws.Range("F513").Value = 18
$ws.Range("G513").Value = 870.48
$ws.Range("B515").Value = 39336.98
$ws.Range("F518").Value = 11
$ws.Range("G518").Value = 2505.58
$ws.Range("F520").Value = 17
$ws.Range("G520").Value = 5857.69
$ws.Range("F521").Value = 571
$ws.Range("G521").Value = 6121.12
$ws.Range("F524").Value = 119
$ws.Range("G524").Value = 7225.68
$ws.Range("F529").Value = 151
$ws.Range("G529").Value = 3699.5
$ws.Range("F531").Value = 44
$ws.Range("G531").Value = 1070.96
$ws.Range("F536").Value = 44
$ws.Range("G536").Value = 25693.8
$ws.Range("B539").Value = 117632.72
$ws.Range("F577").Value = 50
$ws.Range("G577").Value = 3940
$ws.Range("B581").Value = 18326.61
$ws.Range("F586").Value = 50
$ws.Range("G586").Value = 13304.5
$ws.Range("F589").Value = 0
$ws.Range("G589").Value = 0
$ws.Range("B600").Value = 70821.23
$ws.Range("F606").Value = 330
$ws.Range("G606").Value = 1184.7
$ws.Range("F615").Value = 117
$ws.Range("G615").Value = 2797.47
$ws.Range("F618").Value = 440
$ws.Range("G618").Value = 1487.2
$ws.Range("F624").Value = 47
$ws.Range("G624").Value = 4935.94
$ws.Range("B625").Value = 46623.48
$ws.Range("F656").Value = 1
$ws.Range("G656").Value = 27.2
$ws.Range("F657").Value = 27
$ws.Range("G657").Value = 734.4
$ws.Range("F658").Value = 14
$ws.Range("G658").Value = 380.8
$ws.Range("B659").Value = 10268.92
$ws.Range("F719").Value = 52
$ws.Range("G719").Value = 3473.08
$ws.Range("F720").Value = 39
$ws.Range("G720").Value = 3920.67
$ws.Range("B722").Value = 28304.44
$ws.Range("F748").Value = 164
$ws.Range("G748").Value = 21828.4
$ws.Range("F759").Value = 183
$ws.Range("G759").Value = 27656.79
$ws.Range("B767").Value = 570056.14
$ws.Range("F773").Value = 1
$ws.Range("G773").Value = 25.33
$ws.Range("F783").Value = 14
$ws.Range("G783").Value = 2115.82
$ws.Range("F787").Value = 26
$ws.Range("G787").Value = 9823.059999999999
$ws.Range("F788").Value = 333
$ws.Range("G788").Value = 34269.03
$ws.Range("B796").Value = 166677.84
$ws.Range("B855").Value = 5401421.27
$ws.Range("B856").Value = 5401421.27
